$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 26 ("RM 232") entirely; rows below shift up by one.
$ws.Rows.Item(26).Delete()

# After the shift, "SC 92" (originally row 28) is now at row 27. Delete it too.
$ws.Rows.Item(27).Delete()

# Fix up column D values that were re-imputed for the remaining rows.
$ws.Range("D27").Value = -14.6   # SC 101
$ws.Range("D28").Value = ""      # SC 105 (now missing)
$ws.Range("D29").Value = ""      # SC 119 (now missing)
$ws.Range("D30").Value = -13.6   # SC 120
$ws.Range("D32").Value = ""      # SC 193 (now missing)
